$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- source row 27
$ws.Cells.Item(2, 4).Value = 44335
$ws.Cells.Item(2, 8).Value = 'Perfection'
$ws.Cells.Item(2, 10).Value = 100
$ws.Cells.Item(2, 11).Value = 30000
$ws.Cells.Item(2, 12).Value = 32000
$ws.Cells.Item(2, 13).Value = 31000
$ws.Cells.Item(2, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(2, 15).Value = 'Provincia de Huasco'
$ws.Cells.Item(2, 16).Value = 1240

# Row 3 <- source row 9
$ws.Cells.Item(3, 4).Value = 44342
$ws.Cells.Item(3, 8).Value = 'Perfection'
$ws.Cells.Item(3, 10).Value = 60
$ws.Cells.Item(3, 11).Value = 30000
$ws.Cells.Item(3, 12).Value = 32000
$ws.Cells.Item(3, 13).Value = 31000
$ws.Cells.Item(3, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(3, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(3, 16).Value = 1240

# Row 4 <- source row 17
$ws.Cells.Item(4, 4).Value = 44518
$ws.Cells.Item(4, 8).Value = 'Perfection'
$ws.Cells.Item(4, 10).Value = 350
$ws.Cells.Item(4, 11).Value = 14000
$ws.Cells.Item(4, 12).Value = 15000
$ws.Cells.Item(4, 13).Value = 14571
$ws.Cells.Item(4, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(4, 15).Value = 'Región del Maule'
$ws.Cells.Item(4, 16).Value = 583

# Row 5 <- source row 18
$ws.Cells.Item(5, 4).Value = 44519
$ws.Cells.Item(5, 8).Value = 'Perfection'
$ws.Cells.Item(5, 10).Value = 240
$ws.Cells.Item(5, 11).Value = 17000
$ws.Cells.Item(5, 12).Value = 18000
$ws.Cells.Item(5, 13).Value = 17583
$ws.Cells.Item(5, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(5, 15).Value = 'Carahue'
$ws.Cells.Item(5, 16).Value = 703

# Row 6 <- source row 30
$ws.Cells.Item(6, 4).Value = 44659
$ws.Cells.Item(6, 8).Value = 'Sin especificar'
$ws.Cells.Item(6, 10).Value = 140
$ws.Cells.Item(6, 11).Value = 24000
$ws.Cells.Item(6, 12).Value = 25000
$ws.Cells.Item(6, 13).Value = 24571
$ws.Cells.Item(6, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(6, 15).Value = 'Carahue'
$ws.Cells.Item(6, 16).Value = 983

# Row 7 <- source row 28
$ws.Cells.Item(7, 4).Value = 44328
$ws.Cells.Item(7, 8).Value = 'Perfection'
$ws.Cells.Item(7, 10).Value = 100
$ws.Cells.Item(7, 11).Value = 33000
$ws.Cells.Item(7, 12).Value = 34000
$ws.Cells.Item(7, 13).Value = 33500
$ws.Cells.Item(7, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(7, 15).Value = 'Provincia de Huasco'
$ws.Cells.Item(7, 16).Value = 1340

# Row 8 <- source row 13
$ws.Cells.Item(8, 4).Value = 44629
$ws.Cells.Item(8, 8).Value = 'Perfection'
$ws.Cells.Item(8, 10).Value = 35
$ws.Cells.Item(8, 11).Value = 25000
$ws.Cells.Item(8, 12).Value = 26000
$ws.Cells.Item(8, 13).Value = 25429
$ws.Cells.Item(8, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(8, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(8, 16).Value = 1017

# Row 9 <- source row 25
$ws.Cells.Item(9, 4).Value = 44615
$ws.Cells.Item(9, 8).Value = 'Sin especificar'
$ws.Cells.Item(9, 10).Value = 200
$ws.Cells.Item(9, 11).Value = 28000
$ws.Cells.Item(9, 12).Value = 30000
$ws.Cells.Item(9, 13).Value = 29000
$ws.Cells.Item(9, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(9, 15).Value = 'Carahue'
$ws.Cells.Item(9, 16).Value = 1160

# Row 10 <- source row 15
$ws.Cells.Item(10, 4).Value = 44503
$ws.Cells.Item(10, 8).Value = 'Perfection'
$ws.Cells.Item(10, 10).Value = 200
$ws.Cells.Item(10, 11).Value = 15000
$ws.Cells.Item(10, 12).Value = 16000
$ws.Cells.Item(10, 13).Value = 15500
$ws.Cells.Item(10, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(10, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(10, 16).Value = 620

# Row 11 <- source row 19
$ws.Cells.Item(11, 4).Value = 44483
$ws.Cells.Item(11, 8).Value = 'Perfection'
$ws.Cells.Item(11, 10).Value = 220
$ws.Cells.Item(11, 11).Value = 19000
$ws.Cells.Item(11, 12).Value = 20000
$ws.Cells.Item(11, 13).Value = 19455
$ws.Cells.Item(11, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(11, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(11, 16).Value = 778

# Row 12 <- source row 26
$ws.Cells.Item(12, 4).Value = 44545
$ws.Cells.Item(12, 8).Value = 'Perfection'
$ws.Cells.Item(12, 10).Value = 180
$ws.Cells.Item(12, 11).Value = 15000
$ws.Cells.Item(12, 12).Value = 16000
$ws.Cells.Item(12, 13).Value = 15444
$ws.Cells.Item(12, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(12, 15).Value = 'Carahue'
$ws.Cells.Item(12, 16).Value = 618

# Row 13 <- source row 4
$ws.Cells.Item(13, 4).Value = 44595
$ws.Cells.Item(13, 8).Value = 'Perfection'
$ws.Cells.Item(13, 10).Value = 50
$ws.Cells.Item(13, 11).Value = 26000
$ws.Cells.Item(13, 12).Value = 28000
$ws.Cells.Item(13, 13).Value = 27200
$ws.Cells.Item(13, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(13, 15).Value = 'Carahue'
$ws.Cells.Item(13, 16).Value = 1088

# Row 14 <- source row 12
$ws.Cells.Item(14, 4).Value = 44532
$ws.Cells.Item(14, 8).Value = 'Sin especificar'
$ws.Cells.Item(14, 10).Value = 250
$ws.Cells.Item(14, 11).Value = 14000
$ws.Cells.Item(14, 12).Value = 15000
$ws.Cells.Item(14, 13).Value = 14400
$ws.Cells.Item(14, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(14, 15).Value = 'Región del Maule'
$ws.Cells.Item(14, 16).Value = 576

# Row 15 <- source row 20
$ws.Cells.Item(15, 4).Value = 44533
$ws.Cells.Item(15, 8).Value = 'Perfection'
$ws.Cells.Item(15, 10).Value = 80
$ws.Cells.Item(15, 11).Value = 14000
$ws.Cells.Item(15, 12).Value = 15000
$ws.Cells.Item(15, 13).Value = 14375
$ws.Cells.Item(15, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(15, 15).Value = 'Región del Maule'
$ws.Cells.Item(15, 16).Value = 575

# Row 16 <- source row 7
$ws.Cells.Item(16, 4).Value = 44517
$ws.Cells.Item(16, 8).Value = 'Perfection'
$ws.Cells.Item(16, 10).Value = 110
$ws.Cells.Item(16, 11).Value = 17000
$ws.Cells.Item(16, 12).Value = 18000
$ws.Cells.Item(16, 13).Value = 17455
$ws.Cells.Item(16, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(16, 15).Value = 'Región del Maule'
$ws.Cells.Item(16, 16).Value = 698

# Row 17 <- source row 11
$ws.Cells.Item(17, 4).Value = 44657
$ws.Cells.Item(17, 8).Value = 'Sin especificar'
$ws.Cells.Item(17, 10).Value = 250
$ws.Cells.Item(17, 11).Value = 24000
$ws.Cells.Item(17, 12).Value = 25000
$ws.Cells.Item(17, 13).Value = 24400
$ws.Cells.Item(17, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(17, 15).Value = 'Carahue'
$ws.Cells.Item(17, 16).Value = 976

# Row 18 <- source row 2
$ws.Cells.Item(18, 4).Value = 44589
$ws.Cells.Item(18, 8).Value = 'Perfection'
$ws.Cells.Item(18, 10).Value = 160
$ws.Cells.Item(18, 11).Value = 22000
$ws.Cells.Item(18, 12).Value = 23000
$ws.Cells.Item(18, 13).Value = 22500
$ws.Cells.Item(18, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(18, 15).Value = 'Carahue'
$ws.Cells.Item(18, 16).Value = 900

# Row 19 <- source row 5
$ws.Cells.Item(19, 4).Value = 44539
$ws.Cells.Item(19, 8).Value = 'Sin especificar'
$ws.Cells.Item(19, 10).Value = 50
$ws.Cells.Item(19, 11).Value = 13000
$ws.Cells.Item(19, 12).Value = 14000
$ws.Cells.Item(19, 13).Value = 13400
$ws.Cells.Item(19, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(19, 15).Value = 'Región del Maule'
$ws.Cells.Item(19, 16).Value = 536

# Row 20 <- source row 31
$ws.Cells.Item(20, 4).Value = 44162
$ws.Cells.Item(20, 8).Value = 'Sin especificar'
$ws.Cells.Item(20, 10).Value = 100
$ws.Cells.Item(20, 11).Value = 17000
$ws.Cells.Item(20, 12).Value = 18000
$ws.Cells.Item(20, 13).Value = 17500
$ws.Cells.Item(20, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(20, 15).Value = 'Región del Maule'
$ws.Cells.Item(20, 16).Value = 700

# Row 21 <- source row 24
$ws.Cells.Item(21, 4).Value = 44512
$ws.Cells.Item(21, 8).Value = 'Perfection'
$ws.Cells.Item(21, 10).Value = 100
$ws.Cells.Item(21, 11).Value = 14000
$ws.Cells.Item(21, 12).Value = 15000
$ws.Cells.Item(21, 13).Value = 14500
$ws.Cells.Item(21, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(21, 15).Value = 'Región del Maule'
$ws.Cells.Item(21, 16).Value = 580

# Row 22 <- source row 6
$ws.Cells.Item(22, 4).Value = 44540
$ws.Cells.Item(22, 8).Value = 'Sin especificar'
$ws.Cells.Item(22, 10).Value = 110
$ws.Cells.Item(22, 11).Value = 16000
$ws.Cells.Item(22, 12).Value = 17000
$ws.Cells.Item(22, 13).Value = 16545
$ws.Cells.Item(22, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(22, 15).Value = 'Región del Maule'
$ws.Cells.Item(22, 16).Value = 662

# Row 24 <- source row 10
$ws.Cells.Item(24, 4).Value = 44399
$ws.Cells.Item(24, 8).Value = 'Perfection'
$ws.Cells.Item(24, 10).Value = 50
$ws.Cells.Item(24, 11).Value = 39000
$ws.Cells.Item(24, 12).Value = 40000
$ws.Cells.Item(24, 13).Value = 39600
$ws.Cells.Item(24, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(24, 15).Value = 'Provincia de Huasco'
$ws.Cells.Item(24, 16).Value = 1584

# Row 25 <- source row 14
$ws.Cells.Item(25, 4).Value = 44505
$ws.Cells.Item(25, 8).Value = 'Perfection'
$ws.Cells.Item(25, 10).Value = 210
$ws.Cells.Item(25, 11).Value = 6500
$ws.Cells.Item(25, 12).Value = 7000
$ws.Cells.Item(25, 13).Value = 6714
$ws.Cells.Item(25, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(25, 15).Value = 'Región del Maule'
$ws.Cells.Item(25, 16).Value = 269

# Row 26 <- source row 8
$ws.Cells.Item(26, 4).Value = 44673
$ws.Cells.Item(26, 8).Value = 'Sin especificar'
$ws.Cells.Item(26, 10).Value = 220
$ws.Cells.Item(26, 11).Value = 25000
$ws.Cells.Item(26, 12).Value = 26000
$ws.Cells.Item(26, 13).Value = 25455
$ws.Cells.Item(26, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(26, 15).Value = 'Carahue'
$ws.Cells.Item(26, 16).Value = 1018

# Row 27 <- source row 29
$ws.Cells.Item(27, 4).Value = 44643
$ws.Cells.Item(27, 8).Value = 'Perfection'
$ws.Cells.Item(27, 10).Value = 90
$ws.Cells.Item(27, 11).Value = 25000
$ws.Cells.Item(27, 12).Value = 26000
$ws.Cells.Item(27, 13).Value = 25444
$ws.Cells.Item(27, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(27, 15).Value = 'Carahue'
$ws.Cells.Item(27, 16).Value = 1018

# Row 28 <- source row 21
$ws.Cells.Item(28, 4).Value = 44496
$ws.Cells.Item(28, 8).Value = 'Perfection'
$ws.Cells.Item(28, 10).Value = 250
$ws.Cells.Item(28, 11).Value = 14000
$ws.Cells.Item(28, 12).Value = 15000
$ws.Cells.Item(28, 13).Value = 14520
$ws.Cells.Item(28, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(28, 15).Value = 'Provincia de Huasco'
$ws.Cells.Item(28, 16).Value = 581

# Row 29 <- source row 3
$ws.Cells.Item(29, 4).Value = 44671
$ws.Cells.Item(29, 8).Value = 'Perfection'
$ws.Cells.Item(29, 10).Value = 110
$ws.Cells.Item(29, 11).Value = 25000
$ws.Cells.Item(29, 12).Value = 26000
$ws.Cells.Item(29, 13).Value = 25545
$ws.Cells.Item(29, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(29, 15).Value = 'Carahue'
$ws.Cells.Item(29, 16).Value = 1022

# Row 30 <- source row 16
$ws.Cells.Item(30, 4).Value = 44631
$ws.Cells.Item(30, 8).Value = 'Perfection'
$ws.Cells.Item(30, 10).Value = 150
$ws.Cells.Item(30, 11).Value = 24000
$ws.Cells.Item(30, 12).Value = 25000
$ws.Cells.Item(30, 13).Value = 24467
$ws.Cells.Item(30, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(30, 15).Value = 'Carahue'
$ws.Cells.Item(30, 16).Value = 979

# Row 31 <- source row 22
$ws.Cells.Item(31, 4).Value = 44454
$ws.Cells.Item(31, 8).Value = 'Perfection'
$ws.Cells.Item(31, 10).Value = 100
$ws.Cells.Item(31, 11).Value = 36000
$ws.Cells.Item(31, 12).Value = 38000
$ws.Cells.Item(31, 13).Value = 37000
$ws.Cells.Item(31, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(31, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(31, 16).Value = 1480
